$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cell B124: 66,420 -> 66,336 (kept as text, matching existing column formatting) ---
$ws.Cells.Item(124, 2).NumberFormat = "@"
$ws.Cells.Item(124, 2).Value = "66,336"
$ws.Cells.Item(124, 2).Style = "Normal"

# --- Append new row 125: "2025 APR" / "65,214" ---
$ws.Cells.Item(125, 1).NumberFormat = "@"
$ws.Cells.Item(125, 1).Value = "2025 APR"
$ws.Cells.Item(125, 1).Style = "Normal"

$ws.Cells.Item(125, 2).NumberFormat = "@"
$ws.Cells.Item(125, 2).Value = "65,214"
$ws.Cells.Item(125, 2).Style = "Normal"
